$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Value)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '64.257.56'
Set-TextValue $ws.Range('E2') '  -1.92%  '
Set-TextValue $ws.Range('D3') '3.192.60'
Set-TextValue $ws.Range('E3') '  -6.98%  '
Set-TextValue $ws.Range('E4') '  +0.03%  '
Set-TextValue $ws.Range('D5') '558.02'
Set-TextValue $ws.Range('E5') '  -3.95%  '
Set-TextValue $ws.Range('D6') '173.05'
Set-TextValue $ws.Range('E6') '  +0.90%  '
Set-TextValue $ws.Range('E7') '  +0.04%  '
Set-TextValue $ws.Range('D8') '0.605'
Set-TextValue $ws.Range('E8') '  +1.69%  '
Set-TextValue $ws.Range('D9') '3.193.42'
Set-TextValue $ws.Range('E9') '  -6.90%  '
Set-TextValue $ws.Range('E10') '  -4.88%  '
Set-TextValue $ws.Range('D11') '6.61'
Set-TextValue $ws.Range('E11') '  -3.97%  '
Set-TextValue $ws.Range('E12') '  -2.82%  '
Set-TextValue $ws.Range('D13') '3.752.76'
Set-TextValue $ws.Range('E13') '  -6.84%  '
Set-TextValue $ws.Range('E14') '  +0.06%  '
Set-TextValue $ws.Range('D15') '27.35'
Set-TextValue $ws.Range('E15') '  -5.07%  '
Set-TextValue $ws.Range('D16') '64.434.98'
Set-TextValue $ws.Range('E16') '  -1.86%  '
Set-TextValue $ws.Range('E17') '  -4.33%  '
Set-TextValue $ws.Range('D18') '3.208.08'
Set-TextValue $ws.Range('E18') '  -6.57%  '
Set-TextValue $ws.Range('D19') '5.64'
Set-TextValue $ws.Range('E19') '  -4.33%  '
Set-TextValue $ws.Range('D20') '12.99'
Set-TextValue $ws.Range('E20') '  -5.87%  '
Set-TextValue $ws.Range('D21') '354.39'
Set-TextValue $ws.Range('E21') '  -2.86%  '
Set-TextValue $ws.Range('E22') '  -5.63%  '
Set-TextValue $ws.Range('D23') '1.00'
Set-TextValue $ws.Range('E23') '  +0.30%  '
Set-TextValue $ws.Range('D24') '68.75'
Set-TextValue $ws.Range('E24') '  -5.34%  '
Set-TextValue $ws.Range('D25') '0.0000118'
Set-TextValue $ws.Range('E25') '  -2.16%  '
Set-TextValue $ws.Range('D26') '0.502'
Set-TextValue $ws.Range('E26') '  -5.34%  '
Set-TextValue $ws.Range('D27') '9.44'
Set-TextValue $ws.Range('E27') '  -2.96%  '
Set-TextValue $ws.Range('E28') '  -1.57%  '
Set-TextValue $ws.Range('D29') '1.00'
Set-TextValue $ws.Range('E29') '  +0.00%  '
Set-TextValue $ws.Range('B30') 'NEARProtocol'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D30') '5.60'
Set-TextValue $ws.Range('E30') '  -1.31%  '
Set-TextValue $ws.Range('B31') 'USDe'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range('D31') '0.999'
Set-TextValue $ws.Range('E31') '  -0.06%  '
Set-TextValue $ws.Range('E32') '  -4.27%  '
Set-TextValue $ws.Range('D33') '22.13'
Set-TextValue $ws.Range('E33') '  -6.33%  '
Set-TextValue $ws.Range('D34') '6.59'
Set-TextValue $ws.Range('E34') '  -5.93%  '
Set-TextValue $ws.Range('E35') '  -7.53%  '
Set-TextValue $ws.Range('D36') '158.42'
Set-TextValue $ws.Range('E36') '  -1.37%  '
Set-TextValue $ws.Range('D37') '1.43'
Set-TextValue $ws.Range('E37') '  -5.28%  '
Set-TextValue $ws.Range('D38') '0.812'
Set-TextValue $ws.Range('E38') '  -7.64%  '
Set-TextValue $ws.Range('D39') '26.21'
Set-TextValue $ws.Range('E39') '  -9.13%  '
Set-TextValue $ws.Range('D40') '2.51'
Set-TextValue $ws.Range('E40') '  -2.75%  '
Set-TextValue $ws.Range('B41') 'Maker'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D41') '2.653.53'
Set-TextValue $ws.Range('E41') '  -3.59%  '
Set-TextValue $ws.Range('B42') 'Stacks'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D42') '1.67'
Set-TextValue $ws.Range('E42') '  -4.06%  '
Set-TextValue $ws.Range('D43') '6.06'
Set-TextValue $ws.Range('E43') '  -5.25%  '
Set-TextValue $ws.Range('D44') '4.13'
Set-TextValue $ws.Range('E44') '  -6.34%  '
Set-TextValue $ws.Range('B45') 'Hedera'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D45') '0.0653'
Set-TextValue $ws.Range('E45') '  -3.51%  '
Set-TextValue $ws.Range('B46') 'OKB'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D46') '38.92'
Set-TextValue $ws.Range('E46') '  -2.38%  '
Set-TextValue $ws.Range('D47') '320.39'
Set-TextValue $ws.Range('E47') '  -0.47%  '
Set-TextValue $ws.Range('D48') '23.31'
Set-TextValue $ws.Range('E48') '  -2.90%  '
Set-TextValue $ws.Range('D49') '0.0269'
Set-TextValue $ws.Range('E49') '  -6.17%  '
Set-TextValue $ws.Range('D50') '0.102'
Set-TextValue $ws.Range('E50') '  +0.81%  '
Set-TextValue $ws.Range('E51') '  +0.04%  '
